# "various updates including added comments to Build.ps1"
#
# Sheet1 ("My Stuff/Useful git.xlsx") gains three new git-command rows, the
# previously bold/special-font row 2 reverts to plain formatting, and
# column A widens to fit the new (longer) descriptions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("git remote prune origin") loses its custom row height / font
# formatting and goes back to the sheet's plain default style.
$ws.Range("A2:B2").Style = "Normal"
$ws.Rows(2).AutoFit()

# New rows: col A = Description, col B = Command (same layout as existing data).
$ws.Range("A3").Value = "Undo the last commit"
$ws.Range("B3").Value = "git reset head^"

$ws.Range("A4").Value = "Undo the last add"
$ws.Range("B4").Value = "git reset head"

$ws.Range("A5").Value = "Overwrite working directory with what head is currently pointing to"
$ws.Range("B5").Value = "git reset --hard head"

# Column A needs to be wide enough for the longest new description.
$ws.Columns(1).ColumnWidth = 61.5
